$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @("帝国守护之守望苍穹", 45903.3009259259, "3302034a", "907566442"),
    @("无奈",               45903.7280555556, "9b41a46b", "487996763"),
    @("哦",                 45903.8242592593, "fee46bac", "1730863123"),
    @("落叶",               45904.3019791667, "6c3bcd6f", "3382881855"),
    @("我是mc",             45904.5552430556, "46e4be66", "2721614772"),
    @("无执T̶o̶n̶z̶y̶",      45905.6912037037, "72880a09", "3521425739")
)

$r = 29
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
